$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.082.28"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "1.629.52"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.21"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.04"
$ws.Range("E10").Value = "  -0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").Value = "1.857.39"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").Value = "1.624.01"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.67"
$ws.Range("E16").Value = "  -3.54%  "

$ws.Range("D17").Value = "27.059.76"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.70"
$ws.Range("E19").Value = "  -3.32%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"

$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.51"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("E23").Value = "  -1.47%  "

$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.07"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  -1.23%  "

$ws.Range("D34").Value = "1.301.31"
$ws.Range("E34").Value = "  +2.43%  "

$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.542"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.26"
$ws.Range("E41").Value = "  +3.98%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.808"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("E43").Value = "  -2.01%  "

$ws.Range("D44").Value = "1.766.60"
$ws.Range("E44").Value = "  -1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.07"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.37"
$ws.Range("E46").Value = "  -2.29%  "

$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.786"
$ws.Range("E49").Value = "  +17.48%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0513"
$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.54"
$ws.Range("E51").Value = "  -2.14%  "
